# ADD results from server
# Update investment cost values on the "2025", "2030", and "2035" sheets
# with refreshed figures pulled from the server.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 1715.760994200636
$ws.Range("E2").Value = 249348.5034426333
$ws.Range("I2").Value = 138240.852797177
$ws.Range("L2").Value = 448883.1064232461
$ws.Range("M2").Value = 101434.8210067014
$ws.Range("N2").Value = 65159.57982584304
$ws.Range("O2").Value = 62112.65428698476

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 1560.164813693395
$ws.Range("B2").Value = 41737.66617636
$ws.Range("E2").Value = 245152.3900012034
$ws.Range("I2").Value = 225425.3770606828
$ws.Range("L2").Value = 152673.1028312588
$ws.Range("M2").Value = 95084.48612473471
$ws.Range("N2").Value = 30791.72579164192
$ws.Range("O2").Value = 21756.19778172665

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 20479.74258078725
$ws.Range("B2").Value = 14871.78471240001
$ws.Range("E2").Value = 98812.75330607952
$ws.Range("I2").Value = 144259.0521526223
$ws.Range("M2").Value = 54243.78658924496
$ws.Range("N2").Value = 43635.69805396052
$ws.Range("O2").Value = 53147.86046225035
